$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.544408440589905
$ws.Range("B1").Value = 2.42726469039917
$ws.Range("C1").Value = 4.407708168029785
$ws.Range("D1").Value = 1.788809180259705
$ws.Range("E1").Value = 0.8073791861534119
